$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G ("Recorded By") whose two comma-separated names need to
# be swapped in order (per the daily attendance processing commit).
$rows = @(10,11,12,13,14,15,17,18,19,20,21,22,24,26,36,37,38,39,40,41,43,44,45,46,47,48,50,52,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2
    $parts = $val -split ", ", 2
    if ($parts.Length -eq 2) {
        $cell.Value2 = $parts[1] + ", " + $parts[0]
    }
}
